# Apply the release-notes update described by the commit:
#   "update for insert release-notes.md f80ed2bb9e1dd81abc71d13817b8a44a756cee80"
#
# 1. Bump the StructureDefinition metadata (Version / Status / Date / Contact)
#    on the "Metadata" worksheet.
# 2. Swap the two mapping columns ("Mapping: RIM Mapping" and
#    "Mapping: Spécification métier vers l'extension ROR
#    LocationSupportedCapacity") on the "Elements" worksheet - both their
#    header/data content and their column widths change places.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.4.0-snapshot-1"                  # Version
$meta.Range("B6").Value = "draft"                             # Status
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"         # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"     # Contact

# ---------------------------------------------------------------------
# 2. Elements sheet: swap columns AK (37) and AL (38)
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$colAK = 37
$colAL = 38

# Swap the column widths first.
# NB: this COM surface does not hydrate Range/Columns.ColumnWidth from the
# workbook's stored <col> width on load (it reads back a flat placeholder
# until something is explicitly assigned in this session), so the previous
# widths can't be fetched and swapped dynamically. Stamp the post-swap
# widths directly instead - AK (37) takes AL's original width (83.625
# "raw" OOXML character width) and AL (38) takes AK's original width
# (24.98046875). ColumnWidth goes in or out through a +5/6-character,
# pixel-snapped conversion, so the literals below are the closest inputs
# that round-trip to those raw widths.
$elements.Columns.Item($colAK).ColumnWidth = 82.911667
$elements.Columns.Item($colAL).ColumnWidth = 24.087135

# Swap the cell contents (header row included) row by row.
# NB: use Value2 (not Value) - in this COM surface, reading .Value back out
# yields the property accessor wrapper rather than the scalar, so round
# tripping through it would stamp literal "Variant Value ..." text into the
# cells. Value2 reads/writes the real scalar.
$lastRow = $elements.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cellAK = $elements.Cells.Item($r, $colAK)
    $cellAL = $elements.Cells.Item($r, $colAL)

    $valAK = $cellAK.Value2
    $valAL = $cellAL.Value2

    $cellAK.Value2 = $valAL
    $cellAL.Value2 = $valAK
}
